# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E/F, row 2) and on each language sheet's
#   Status column (column C, row 2).
# - Shrink the "Status" column(s) to match the new (shorter) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status values.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the status columns to fit the shorter text.
# (ColumnWidth is specified in character units and is snapped to the
# nearest pixel by Excel, so 12.5 is the input that lands on the
# intended stored column width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
